# edit.ps1 - scheduled-runner refresh of market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ] -> H:N)
# across the per-job Leve tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Only H:N are touched; A:G (leve name/item/level/exp/gil/amount/item id) are
# left untouched since the upstream market-board pull only refreshed pricing.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 3377.889
$ws.Range("I40").Value = 3377.889
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3377.889
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3202.889
$ws.Range("N40").ClearContents()
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 7900
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 7800
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 7800
$ws.Range("M43").Value = -7931
$ws.Range("N43").Value = -7938
# Row 53: No Accounting for Waste
$ws.Range("H53").Value = 298.5
$ws.Range("I53").Value = 198.5
$ws.Range("J53").Value = 398.5
$ws.Range("K53").Value = 198.5
$ws.Range("L53").Value = 398.5
$ws.Range("M53").Value = 438.5
$ws.Range("N53").Value = -1672.5
# Row 82: Rolling on Initiative
$ws.Range("H82").Value = 5522.933
$ws.Range("I82").Value = 4218.846
$ws.Range("J82").Value = 13999.5
$ws.Range("K82").Value = 12656.538
$ws.Range("L82").Value = 41998.5
$ws.Range("M82").Value = -12250.538
$ws.Range("N82").Value = -42810.5
# Row 85: Darkly Dreaming Dexterity (L)
$ws.Range("H85").Value = 5522.933
$ws.Range("I85").Value = 4218.846
$ws.Range("J85").Value = 13999.5
$ws.Range("K85").Value = 12656.538
$ws.Range("L85").Value = 41998.5
$ws.Range("M85").Value = -11252.538
$ws.Range("N85").Value = -44806.5
# Row 98: The Dotted Line
$ws.Range("H98").Value = 1331.238
$ws.Range("I98").Value = 1197.5
$ws.Range("J98").Value = 4006
$ws.Range("K98").Value = 1197.5
$ws.Range("L98").Value = 4006
$ws.Range("M98").Value = 300.5
$ws.Range("N98").Value = -7002
# Row 100: Asking for a Friend
$ws.Range("H100").Value = 710
$ws.Range("I100").Value = 450
$ws.Range("K100").Value = 450
$ws.Range("M100").Value = 91
# Row 122: Wishful Inking
$ws.Range("H122").Value = 1331.238
$ws.Range("I122").Value = 1197.5
$ws.Range("J122").Value = 4006
$ws.Range("K122").Value = 3592.5
$ws.Range("L122").Value = 12018
$ws.Range("M122").Value = -1142.5
$ws.Range("N122").Value = -16918
# Row 129: Practical Command
$ws.Range("H129").Value = 1340.6364
$ws.Range("I129").Value = 694.3333
$ws.Range("K129").Value = 2082.9999
$ws.Range("M129").Value = 2917.0001
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2273.25
$ws.Range("I132").Value = 2273.25
$ws.Range("K132").Value = 6819.75
$ws.Range("M132").Value = -4289.75
# Row 138: All-night Crafting
$ws.Range("H138").Value = 1863.6666
$ws.Range("I138").Value = 1863.6666
$ws.Range("K138").Value = 5590.9998
$ws.Range("M138").Value = -450.9997999999996
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 759.6667
$ws.Range("I141").Value = 759.6667
$ws.Range("K141").Value = 2279.0001
$ws.Range("M141").Value = 2900.9999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 97: Ore for Me
$ws.Range("H97").Value = 286.58334
$ws.Range("I97").Value = 323.55554
$ws.Range("J97").Value = 175.66667
$ws.Range("K97").Value = 323.55554
$ws.Range("L97").Value = 175.66667
$ws.Range("M97").Value = 172.44446
$ws.Range("N97").Value = -1167.66667

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 3035.5173
$ws.Range("I20").Value = 3053.1
$ws.Range("J20").Value = 2996.4443
$ws.Range("K20").Value = 3053.1
$ws.Range("L20").Value = 2996.4443
$ws.Range("M20").Value = -2806.1
$ws.Range("N20").Value = -3490.4443
# Row 94: High Steal
$ws.Range("H94").Value = 26789.39
$ws.Range("I94").Value = 5007.1514
$ws.Range("K94").Value = 5007.1514
$ws.Range("M94").Value = -4556.1514
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3055
$ws.Range("I105").Value = 2739.5
$ws.Range("K105").Value = 2739.5
$ws.Range("M105").Value = -992.5
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 32190230
$ws.Range("I134").Value = 39617950
$ws.Range("J134").Value = 3431.3333
$ws.Range("K134").Value = 118853850
$ws.Range("L134").Value = 10293.9999
$ws.Range("M134").Value = -118851315
$ws.Range("N134").Value = -15363.9999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 8254.174000000001
$ws.Range("I31").Value = 6495.143
$ws.Range("J31").Value = 9023.75
$ws.Range("K31").Value = 6495.143
$ws.Range("L31").Value = 9023.75
$ws.Range("M31").Value = -6200.143
$ws.Range("N31").Value = -9613.75
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 8254.174000000001
$ws.Range("I34").Value = 6495.143
$ws.Range("J34").Value = 9023.75
$ws.Range("K34").Value = 6495.143
$ws.Range("L34").Value = 9023.75
$ws.Range("M34").Value = -6293.143
$ws.Range("N34").Value = -9427.75
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 23814830
$ws.Range("I58").Value = 25005556
$ws.Range("K58").Value = 25005556
$ws.Range("M58").Value = -25005353
# Row 136: Turali Quality
$ws.Range("H136").Value = 23814830
$ws.Range("I136").Value = 25005556
$ws.Range("K136").Value = 75016668
$ws.Range("M136").Value = -75014118

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 58595.89
$ws.Range("I5").Value = 100627.6
$ws.Range("J5").Value = 6056.25
$ws.Range("K5").Value = 301882.8
$ws.Range("L5").Value = 18168.75
$ws.Range("M5").Value = -301770.8
$ws.Range("N5").Value = -18392.75
# Row 56: Culture Club
$ws.Range("H56").Value = 242380.56
$ws.Range("I56").Value = 242380.56
$ws.Range("K56").Value = 242380.56
$ws.Range("M56").Value = -241850.56
# Row 64: The Aroma of Faith
$ws.Range("H64").Value = 917670.4399999999
$ws.Range("J64").Value = 10799.8
$ws.Range("L64").Value = 32399.4
$ws.Range("N64").Value = -32939.39999999999
# Row 67: Soup's On (L)
$ws.Range("H67").Value = 917670.4399999999
$ws.Range("J67").Value = 10799.8
$ws.Range("L67").Value = 32399.4
$ws.Range("N67").Value = -34271.39999999999
# Row 108: Meet for Meat
$ws.Range("H108").Value = 676.2857
$ws.Range("I108").Value = 676.2857
$ws.Range("K108").Value = 2028.8571
$ws.Range("M108").Value = 851.1428999999998
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 58595.89
$ws.Range("I135").Value = 100627.6
$ws.Range("J135").Value = 6056.25
$ws.Range("K135").Value = 905648.4
$ws.Range("L135").Value = 54506.25
$ws.Range("M135").Value = -903113.4
$ws.Range("N135").Value = -59576.25
# Row 137: Creative Chocolate
$ws.Range("H137").Value = 20001720
$ws.Range("I137").Value = 33334666
$ws.Range("J137").Value = 2299
$ws.Range("K137").Value = 100003998
$ws.Range("L137").Value = 6897
$ws.Range("M137").Value = -99998898
$ws.Range("N137").Value = -17097

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 2199.4546
$ws.Range("I97").Value = 1855.5
$ws.Range("K97").Value = 1855.5
$ws.Range("M97").Value = -1359.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 10000
$ws.Range("K82").Value = 10000
$ws.Range("M82").Value = -9639
# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 10000
$ws.Range("K85").Value = 10000
$ws.Range("M85").Value = -8752
# Row 117: If I Could Walk a Thousand Malms
$ws.Range("H117").Value = 99392
$ws.Range("J117").Value = 99392
$ws.Range("L117").Value = 99392
$ws.Range("N117").Value = -108570
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 4905812.5
$ws.Range("I132").Value = 5003828.5
$ws.Range("K132").Value = 15011485.5
$ws.Range("M132").Value = -15008955.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 41: Half Is the New Double
$ws.Range("H41").Value = 48799.8
$ws.Range("J41").Value = 48799.8
$ws.Range("L41").Value = 48799.8
$ws.Range("N41").Value = -49579.8
# Row 45: Private Concerns
$ws.Range("H45").Value = 20000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
